$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: estado Falta -> Completado (copy the "Completado" look from H3) + fecha de estado ---
$ws.Range("H3").Copy()
$ws.Range("H5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H5").Value = "Completado"
$ws.Range("I5").Value = 44140

# --- Row 7: estado Falta -> Completado + fecha de estado ---
$ws.Range("H3").Copy()
$ws.Range("H7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H7").Value = "Completado"
$ws.Range("I7").Value = 44143

# --- Rows whose "Fecha de estado" (column I) moved to 44143 ---
$ws.Range("I8").Value = 44143
$ws.Range("I9").Value = 44143
$ws.Range("I10").Value = 44143
$ws.Range("I11").Value = 44143
$ws.Range("I12").Value = 44143
$ws.Range("I13").Value = 44143
$ws.Range("I14").Value = 44143
$ws.Range("I15").Value = 44143

$ws.Range("I30").Value = 44143

$ws.Range("I35").Value = 44143
$ws.Range("I36").Value = 44143
$ws.Range("I37").Value = 44143
$ws.Range("I38").Value = 44143
$ws.Range("I39").Value = 44143

# Clear the clipboard marquee / leftover copy state
$excel.CutCopyMode = $false

# Match the author's final cursor position/selection
$ws.Range("I39").Select()
